# Apply the commit "Slides for tokens and react":
#  1. Bump the cached date-field text in the Slide Master and every
#     Slide Layout from 27/09/2023 -> 28/09/2023.
#  2. Re-title the three "Create a User API Controller continued" slides
#     to distinguish Admin / Login / Register variants.

$p = $ppt.ActivePresentation

$newDate = "28/09/2023"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        if ($shp.Type -eq 14) {
            try {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePh = $true
                }
            } catch {
                $isDatePh = $false
            }
        }
        if ($isDatePh -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# 1a. Slide Master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# 1b. Every Slide Layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# 2. Retitle the three "Create a User API Controller continued" slides.
$titleUpdates = @{
    10 = "Create a User API Controller - Admin"
    11 = "Create a User API Controller - Login"
    12 = "Create a User API Controller - Register"
}

foreach ($slideIdx in $titleUpdates.Keys) {
    $slide = $p.Slides.Item($slideIdx)
    $shp = $slide.Shapes.Item(1)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "Create a User API Controller continued") {
        $shp.TextFrame.TextRange.Text = $titleUpdates[$slideIdx]
    }
}
